# Auto-generated edit script applying the Goblin_Profits.xlsx diff.
# For each affected cell: plain value changes use .Value assignment;
# the one cell removed entirely in the diff (CUL!N113) uses ClearContents();
# the newly-introduced cells (BSM!N6, WVR!M74, WVR!M77, WVR!M94) are created
# via plain .Value assignment since the range doesn't exist yet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 200004400
$ws.Range("I69").Value = 3997.5
$ws.Range("K69").Value = 11992.5
$ws.Range("M69").Value = -11118.5
$ws.Range("H72").Value = 200004400
$ws.Range("I72").Value = 3997.5
$ws.Range("K72").Value = 35977.5
$ws.Range("M72").Value = -31609.5
$ws.Range("H74").Value = 3739.2
$ws.Range("I74").Value = 3549
$ws.Range("K74").Value = 3549
$ws.Range("M74").Value = -2613
$ws.Range("H77").Value = 3739.2
$ws.Range("I77").Value = 3549
$ws.Range("K77").Value = 17745
$ws.Range("M77").Value = -13065
$ws.Range("H80").Value = 2241.5
$ws.Range("I80").Value = 605.8182
$ws.Range("J80").Value = 3299.8823
$ws.Range("K80").Value = 1817.4546
$ws.Range("L80").Value = 9899.6469
$ws.Range("M80").Value = -819.4546
$ws.Range("N80").Value = -11895.6469
$ws.Range("H83").Value = 2241.5
$ws.Range("I83").Value = 605.8182
$ws.Range("J83").Value = 3299.8823
$ws.Range("K83").Value = 5452.3638
$ws.Range("L83").Value = 29698.9407
$ws.Range("M83").Value = -460.3638000000001
$ws.Range("N83").Value = -39682.94070000001
$ws.Range("H132").Value = 4587
$ws.Range("I132").Value = 6299.75
$ws.Range("K132").Value = 18899.25
$ws.Range("M132").Value = -16369.25
$ws.Range("H135").Value = 3661
$ws.Range("J135").Value = 3990
$ws.Range("L135").Value = 35910
$ws.Range("N135").Value = -40980

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 268
$ws.Range("I5").Value = 212.5
$ws.Range("K5").Value = 212.5
$ws.Range("M5").Value = -100.5
$ws.Range("H32").Value = 4176.9307
$ws.Range("I32").Value = 3198.913
$ws.Range("K32").Value = 3198.913
$ws.Range("M32").Value = -2911.913
$ws.Range("H45").Value = 15578.875
$ws.Range("I45").Value = 26657.75
$ws.Range("K45").Value = 26657.75
$ws.Range("M45").Value = -26280.75
$ws.Range("H95").Value = 73818.60000000001
$ws.Range("J95").Value = 73818.60000000001
$ws.Range("L95").Value = 73818.60000000001
$ws.Range("N95").Value = -79310.60000000001
$ws.Range("H110").Value = 2569.8696
$ws.Range("I110").Value = 2900.7896
$ws.Range("J110").Value = 998
$ws.Range("K110").Value = 2900.7896
$ws.Range("L110").Value = 998
$ws.Range("M110").Value = -855.7896000000001
$ws.Range("N110").Value = -5088
$ws.Range("H132").Value = 7105.4688
$ws.Range("I132").Value = 6108.963
$ws.Range("K132").Value = 18326.889
$ws.Range("M132").Value = -15796.889

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 268
$ws.Range("I4").Value = 212.5
$ws.Range("K4").Value = 212.5
$ws.Range("M4").Value = -97.5
$ws.Range("H6").Value = 10000
$ws.Range("J6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("N6").Value = -10226
$ws.Range("H13").Value = 49945
$ws.Range("J13").Value = 49945
$ws.Range("L13").Value = 49945
$ws.Range("N13").Value = -50281
$ws.Range("H57").Value = 82000
$ws.Range("J57").Value = 82000
$ws.Range("L57").Value = 82000
$ws.Range("N57").Value = -83440
$ws.Range("H60").Value = 41500
$ws.Range("J60").Value = 41500
$ws.Range("L60").Value = 41500
$ws.Range("N60").Value = -42698
$ws.Range("H105").Value = 6608
$ws.Range("I105").Value = 7855.2
$ws.Range("J105").Value = 3490
$ws.Range("K105").Value = 7855.2
$ws.Range("L105").Value = 3490
$ws.Range("M105").Value = -6108.2
$ws.Range("N105").Value = -6984
$ws.Range("H136").Value = 82000
$ws.Range("J136").Value = 82000
$ws.Range("L136").Value = 82000
$ws.Range("N136").Value = -92200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 39996.5
$ws.Range("J52").Value = 39996.5
$ws.Range("L52").Value = 39996.5
$ws.Range("N52").Value = -40584.5
$ws.Range("H132").Value = 3187.8333
$ws.Range("I132").Value = 3187.8333
$ws.Range("K132").Value = 9563.499899999999
$ws.Range("M132").Value = -7033.499899999999
$ws.Range("H134").Value = 3700.0625
$ws.Range("I134").Value = 2725.2
$ws.Range("K134").Value = 8175.599999999999
$ws.Range("M134").Value = -5640.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1172866.8
$ws.Range("I4").Value = 18929.553
$ws.Range("J4").Value = 22866886
$ws.Range("K4").Value = 56788.659
$ws.Range("L4").Value = 68600658
$ws.Range("M4").Value = -56676.659
$ws.Range("N4").Value = -68600882
$ws.Range("H6").Value = 146.5
$ws.Range("I6").Value = 153.14285
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 459.42855
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = -346.42855
$ws.Range("N6").Value = -526
$ws.Range("H10").Value = 270
$ws.Range("I10").Value = 292
$ws.Range("K10").Value = 876
$ws.Range("M10").Value = -737
$ws.Range("H45").Value = 6418.3335
$ws.Range("J45").Value = 8516.5
$ws.Range("L45").Value = 25549.5
$ws.Range("N45").Value = -26613.5
$ws.Range("H113").Value = 6000
$ws.Range("I113").Value = 6000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 18000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -15830
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1774.4117
$ws.Range("I97").Value = 1774.4117
$ws.Range("K97").Value = 1774.4117
$ws.Range("M97").Value = -1278.4117
$ws.Range("H102").Value = 2388.6667
$ws.Range("I102").Value = 1687.25
$ws.Range("K102").Value = 1687.25
$ws.Range("M102").Value = -65.25
$ws.Range("H126").Value = 2799.1428
$ws.Range("I126").Value = 2600
$ws.Range("K126").Value = 7800
$ws.Range("M126").Value = -5330
$ws.Range("H132").Value = 9137.25
$ws.Range("I132").Value = 5516.5
$ws.Range("K132").Value = 16549.5
$ws.Range("M132").Value = -14019.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 44228.332
$ws.Range("J76").Value = 44228.332
$ws.Range("L76").Value = 44228.332
$ws.Range("N76").Value = -44904.332
$ws.Range("H79").Value = 44228.332
$ws.Range("J79").Value = 44228.332
$ws.Range("L79").Value = 44228.332
$ws.Range("N79").Value = -46568.332
$ws.Range("H100").Value = 8412.375
$ws.Range("I100").Value = 7699.75
$ws.Range("J100").Value = 9125
$ws.Range("K100").Value = 7699.75
$ws.Range("L100").Value = 9125
$ws.Range("M100").Value = -7158.75
$ws.Range("N100").Value = -10207
$ws.Range("H127").Value = 121499.5
$ws.Range("J127").Value = 121499.5
$ws.Range("L127").Value = 121499.5
$ws.Range("N127").Value = -131419.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10812.324
$ws.Range("I62").Value = 3449.8333
$ws.Range("J62").Value = 12237.322
$ws.Range("K62").Value = 3449.8333
$ws.Range("L62").Value = 12237.322
$ws.Range("M62").Value = -2825.8333
$ws.Range("N62").Value = -13485.322
$ws.Range("H65").Value = 10812.324
$ws.Range("I65").Value = 3449.8333
$ws.Range("J65").Value = 12237.322
$ws.Range("K65").Value = 17249.1665
$ws.Range("L65").Value = 61186.61
$ws.Range("M65").Value = -14129.1665
$ws.Range("N65").Value = -67426.61
$ws.Range("H69").Value = 29972.666
$ws.Range("J69").Value = 29972.666
$ws.Range("L69").Value = 29972.666
$ws.Range("N69").Value = -31470.666
$ws.Range("H72").Value = 29972.666
$ws.Range("J72").Value = 29972.666
$ws.Range("L72").Value = 89917.99800000001
$ws.Range("N72").Value = -97405.99800000001
$ws.Range("H74").Value = 22798.889
$ws.Range("I74").Value = 17998
$ws.Range("J74").Value = 23399
$ws.Range("K74").Value = 17998
$ws.Range("L74").Value = 23399
$ws.Range("M74").Value = -17062
$ws.Range("N74").Value = -25271
$ws.Range("H77").Value = 22798.889
$ws.Range("I77").Value = 17998
$ws.Range("J77").Value = 23399
$ws.Range("K77").Value = 53994
$ws.Range("L77").Value = 70197
$ws.Range("M77").Value = -49314
$ws.Range("N77").Value = -79557
$ws.Range("H81").Value = 3096.8
$ws.Range("I81").Value = 3065.8333
$ws.Range("K81").Value = 6131.6666
$ws.Range("M81").Value = -5070.6666
$ws.Range("H84").Value = 3096.8
$ws.Range("I84").Value = 3065.8333
$ws.Range("K84").Value = 30658.333
$ws.Range("M84").Value = -25354.333
$ws.Range("H94").Value = 56640
$ws.Range("I94").Value = 29300
$ws.Range("J94").Value = 61196.668
$ws.Range("K94").Value = 29300
$ws.Range("L94").Value = 61196.668
$ws.Range("M94").Value = -28399
$ws.Range("N94").Value = -62998.668
$ws.Range("H100").Value = 1311
$ws.Range("I100").Value = 766.6667
$ws.Range("K100").Value = 1533.3334
$ws.Range("M100").Value = -992.3334
$ws.Range("H101").Value = 18981.4
$ws.Range("J101").Value = 18981.4
$ws.Range("L101").Value = 18981.4
$ws.Range("N101").Value = -25471.4
$ws.Range("H133").Value = 127000
$ws.Range("J133").Value = 127000
$ws.Range("L133").Value = 127000
$ws.Range("N133").Value = -137120

